# Atualização automática da planilha
#
# On the "Organograma" sheet:
# 1) Fill in the (previously blank) "Time Projeto" / "Key User N1" columns
#    (A and B) for rows 36-42, mirroring the formatting/content already
#    used above them in row 35, and move the active selection to A36:A42.
# 2) Rename the last indicator label in column E (E42) from "Custos" to
#    "Planejamento e Custos".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organograma")

for ($r = 36; $r -le 42; $r++) {
    $ws.Cells.Item(35, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = "Time Projeto"

    $ws.Cells.Item($r, 2).Value = "Key User N1"
}

$ws.Range("A36:A42").Select()

$ws.Cells.Item(42, 5).Value = "Planejamento e Custos"
